{"js": "// Newly downloaded tc/tcn/tl ids were renamed (the stray \"a\" before the\n// trailing digit was dropped): the \"<id>...</id>\" markup for each record\n// collapses from three separately-formatted runs (\"<id>\", the bare id\n// value, \"</id>\") into a single run carrying the id-run's Courier New /\n// 7f6000 formatting and the full updated \"<id>value</id>\" text.\n//\n// Strategy: search the body for the whole visible \"<id>OLD</id>\" span\n// (the search matches across run boundaries) and replace that found\n// range's text in one shot via insertText(..., replace). Word merges a\n// multi-run range like this into a single run that takes on the range's\n// leading (first) run formatting - exactly the merge the diff shows.\n\nasync function replaceIdSpan(oldId, newId) {\n  const oldSpan = `<id>${oldId}</id>`;\n  const newSpan = `<id>${newId}</id>`;\n\n  const results = context.document.body.search(oldSpan, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length > 0) {\n    results.items[0].insertText(newSpan, Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n\nawait replaceIdSpan(\"p085v_a4\", \"p085v_4\");\nawait replaceIdSpan(\"p086v_a1\", \"p086v_1\");\n", "ps1": "# Newly downloaded tc/tcn/tl ids were renamed (the stray \"a\" before the\n# trailing digit was dropped): the <id>...</id> markup for each record\n# collapses from three separately-formatted runs (\"<id>\", the bare id\n# value, \"</id>\") into a single run carrying the id-run's Courier New /\n# 7f6000 formatting and the full updated \"<id>value</id>\" text.\n#\n# Strategy: find the whole visible \"<id>OLD</id>\" span (it reads across\n# the run boundaries) and overwrite that range's Text in one shot. Word\n# collapses a multi-run range assignment like this into a single run that\n# takes on the range's leading (first) run formatting - exactly the\n# merge the diff shows.\n\n$d = $word.ActiveDocument\n\nfunction Replace-IdSpan($doc, [string]$oldId, [string]$newId) {\n    $oldSpan = \"<id>\" + $oldId + \"</id>\"\n    $newSpan = \"<id>\" + $newId + \"</id>\"\n\n    $r = $doc.Content\n    $find = $r.Find\n    $find.ClearFormatting()\n    $find.MatchCase = $true\n    $find.Text = $oldSpan\n    $found = $find.Execute()\n\n    if ($found) {\n        $r.Text = $newSpan\n    }\n}\n\nReplace-IdSpan $d \"p085v_a4\" \"p085v_4\"\nReplace-IdSpan $d \"p086v_a1\" \"p086v_1\"\n"}
